$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1439.3
$ws.Cells.Item(137, 9).Value = 1204.3684
$ws.Cells.Item(137, 11).Value = 3613.1052
$ws.Cells.Item(137, 13).Value = -1063.1052

$ws.Cells.Item(138, 8).Value = 3586.647
$ws.Cells.Item(138, 9).Value = 793.75
$ws.Cells.Item(138, 10).Value = 7209.8647
$ws.Cells.Item(138, 11).Value = 2381.25
$ws.Cells.Item(138, 12).Value = 21629.5941
$ws.Cells.Item(138, 13).Value = 2758.75
$ws.Cells.Item(138, 14).Value = -31909.5941

$ws.Cells.Item(141, 8).Value = 2356.5403
$ws.Cells.Item(141, 9).Value = 1010.8939
$ws.Cells.Item(141, 11).Value = 3032.6817
$ws.Cells.Item(141, 13).Value = 2147.3183

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7246.8706
$ws.Cells.Item(32, 9).Value = 3192.4507
$ws.Cells.Item(32, 11).Value = 3192.4507
$ws.Cells.Item(32, 13).Value = -2905.4507

$ws.Cells.Item(63, 8).Value = 2342.3484
$ws.Cells.Item(63, 9).Value = 2332.1775
$ws.Cells.Item(63, 10).Value = 2500
$ws.Cells.Item(63, 11).Value = 2332.1775
$ws.Cells.Item(63, 12).Value = 2500
$ws.Cells.Item(63, 13).Value = -1646.1775
$ws.Cells.Item(63, 14).Value = -3872

$ws.Cells.Item(66, 8).Value = 2342.3484
$ws.Cells.Item(66, 9).Value = 2332.1775
$ws.Cells.Item(66, 10).Value = 2500
$ws.Cells.Item(66, 11).Value = 11660.8875
$ws.Cells.Item(66, 12).Value = 12500
$ws.Cells.Item(66, 13).Value = -8228.887499999999
$ws.Cells.Item(66, 14).Value = -19364

$ws.Cells.Item(74, 8).Value = 811
$ws.Cells.Item(74, 9).Value = 684.59375
$ws.Cells.Item(74, 10).Value = 1620
$ws.Cells.Item(74, 11).Value = 684.59375
$ws.Cells.Item(74, 12).Value = 1620
$ws.Cells.Item(74, 13).Value = 189.40625
$ws.Cells.Item(74, 14).Value = -3368

$ws.Cells.Item(77, 8).Value = 811
$ws.Cells.Item(77, 9).Value = 684.59375
$ws.Cells.Item(77, 10).Value = 1620
$ws.Cells.Item(77, 11).Value = 3422.96875
$ws.Cells.Item(77, 12).Value = 8100
$ws.Cells.Item(77, 13).Value = 945.03125
$ws.Cells.Item(77, 14).Value = -16836

$ws.Cells.Item(88, 8).Value = 2700
$ws.Cells.Item(88, 9).Value = 1300
$ws.Cells.Item(88, 10).Value = 3209.0908
$ws.Cells.Item(88, 11).Value = 1300
$ws.Cells.Item(88, 12).Value = 3209.0908
$ws.Cells.Item(88, 13).Value = -894
$ws.Cells.Item(88, 14).Value = -4021.0908

$ws.Cells.Item(91, 8).Value = 2700
$ws.Cells.Item(91, 9).Value = 1300
$ws.Cells.Item(91, 10).Value = 3209.0908
$ws.Cells.Item(91, 11).Value = 1300
$ws.Cells.Item(91, 12).Value = 3209.0908
$ws.Cells.Item(91, 13).Value = 104
$ws.Cells.Item(91, 14).Value = -6017.0908

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2457.375
$ws.Cells.Item(99, 9).Value = 2610
$ws.Cells.Item(99, 10).Value = 1999.5
$ws.Cells.Item(99, 11).Value = 2610
$ws.Cells.Item(99, 12).Value = 1999.5
$ws.Cells.Item(99, 13).Value = -1112
$ws.Cells.Item(99, 14).Value = -4995.5

$ws.Cells.Item(107, 8).Value = 1127.4
$ws.Cells.Item(107, 9).Value = 1063.875
$ws.Cells.Item(107, 10).Value = 1381.5
$ws.Cells.Item(107, 11).Value = 1063.875
$ws.Cells.Item(107, 12).Value = 1381.5
$ws.Cells.Item(107, 13).Value = 856.125
$ws.Cells.Item(107, 14).Value = -5221.5

$ws.Cells.Item(134, 8).Value = 1561.0952
$ws.Cells.Item(134, 9).Value = 1181.2084
$ws.Cells.Item(134, 10).Value = 2067.611
$ws.Cells.Item(134, 11).Value = 3543.6252
$ws.Cells.Item(134, 12).Value = 6202.833
$ws.Cells.Item(134, 13).Value = -1008.6252
$ws.Cells.Item(134, 14).Value = -11272.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2312.578
$ws.Cells.Item(31, 9).Value = 1192.25
$ws.Cells.Item(31, 10).Value = 5070.3076
$ws.Cells.Item(31, 11).Value = 1192.25
$ws.Cells.Item(31, 12).Value = 5070.3076
$ws.Cells.Item(31, 13).Value = -897.25
$ws.Cells.Item(31, 14).Value = -5660.3076

$ws.Cells.Item(34, 8).Value = 2312.578
$ws.Cells.Item(34, 9).Value = 1192.25
$ws.Cells.Item(34, 10).Value = 5070.3076
$ws.Cells.Item(34, 11).Value = 1192.25
$ws.Cells.Item(34, 12).Value = 5070.3076
$ws.Cells.Item(34, 13).Value = -990.25
$ws.Cells.Item(34, 14).Value = -5474.3076

$ws.Cells.Item(132, 8).Value = 1468.6285
$ws.Cells.Item(132, 9).Value = 1061.1482
$ws.Cells.Item(132, 10).Value = 2843.875
$ws.Cells.Item(132, 11).Value = 3183.4446
$ws.Cells.Item(132, 12).Value = 8531.625
$ws.Cells.Item(132, 13).Value = -653.4446000000003
$ws.Cells.Item(132, 14).Value = -13591.625

$ws.Cells.Item(134, 8).Value = 1163.3617
$ws.Cells.Item(134, 9).Value = 1020.62164
$ws.Cells.Item(134, 10).Value = 1691.5
$ws.Cells.Item(134, 11).Value = 3061.86492
$ws.Cells.Item(134, 12).Value = 5074.5
$ws.Cells.Item(134, 13).Value = -526.86492
$ws.Cells.Item(134, 14).Value = -10144.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 886.32434
$ws.Cells.Item(5, 9).Value = 322.2258
$ws.Cells.Item(5, 10).Value = 3800.8333
$ws.Cells.Item(5, 11).Value = 966.6774
$ws.Cells.Item(5, 12).Value = 11402.4999
$ws.Cells.Item(5, 13).Value = -854.6774
$ws.Cells.Item(5, 14).Value = -11626.4999

$ws.Cells.Item(131, 8).Value = 890.6667
$ws.Cells.Item(131, 9).Value = 525.75
$ws.Cells.Item(131, 10).Value = 994.9286
$ws.Cells.Item(131, 11).Value = 1577.25
$ws.Cells.Item(131, 12).Value = 2984.7858
$ws.Cells.Item(131, 13).Value = 3462.75
$ws.Cells.Item(131, 14).Value = -13064.7858

$ws.Cells.Item(135, 8).Value = 886.32434
$ws.Cells.Item(135, 9).Value = 322.2258
$ws.Cells.Item(135, 10).Value = 3800.8333
$ws.Cells.Item(135, 11).Value = 2900.0322
$ws.Cells.Item(135, 12).Value = 34207.4997
$ws.Cells.Item(135, 13).Value = -365.0322000000001
$ws.Cells.Item(135, 14).Value = -39277.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 5389.9287
$ws.Cells.Item(132, 9).Value = 6547.1577
$ws.Cells.Item(132, 10).Value = 2946.889
$ws.Cells.Item(132, 11).Value = 19641.4731
$ws.Cells.Item(132, 12).Value = 8840.667000000001
$ws.Cells.Item(132, 13).Value = -17111.4731
$ws.Cells.Item(132, 14).Value = -13900.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 1773.1904
$ws.Cells.Item(136, 9).Value = 889.7931
$ws.Cells.Item(136, 10).Value = 3743.8462
$ws.Cells.Item(136, 11).Value = 2669.3793
$ws.Cells.Item(136, 12).Value = 11231.5386
$ws.Cells.Item(136, 13).Value = -119.3793000000001
$ws.Cells.Item(136, 14).Value = -16331.5386

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 6015.077
$ws.Cells.Item(81, 9).Value = 11081.637
$ws.Cells.Item(81, 10).Value = 2299.6
$ws.Cells.Item(81, 11).Value = 22163.274
$ws.Cells.Item(81, 12).Value = 4599.2
$ws.Cells.Item(81, 13).Value = -21102.274
$ws.Cells.Item(81, 14).Value = -6721.2

$ws.Cells.Item(84, 8).Value = 6015.077
$ws.Cells.Item(84, 9).Value = 11081.637
$ws.Cells.Item(84, 10).Value = 2299.6
$ws.Cells.Item(84, 11).Value = 110816.37
$ws.Cells.Item(84, 12).Value = 22996
$ws.Cells.Item(84, 13).Value = -105512.37
$ws.Cells.Item(84, 14).Value = -33604

$ws.Cells.Item(136, 8).Value = 2600.4033
$ws.Cells.Item(136, 9).Value = 789.0222
$ws.Cells.Item(136, 10).Value = 7395.2354
$ws.Cells.Item(136, 11).Value = 2367.0666
$ws.Cells.Item(136, 12).Value = 22185.7062
$ws.Cells.Item(136, 13).Value = 182.9333999999999
$ws.Cells.Item(136, 14).Value = -27285.7062
